$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cryptocurrency price/volume data per upstream scrape
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.752.01"
$ws.Range("E2").Value = "  -1.82%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.767.59"
$ws.Range("E3").Value = "  -1.68%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.019"
$ws.Range("E4").Value = "  +1.62%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "327.17"
$ws.Range("E5").Value = "  -3.48%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.012"
$ws.Range("E6").Value = "  +1.24%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4311"
$ws.Range("E7").Value = "  -5.55%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3669"
$ws.Range("E8").Value = "  +1.97%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "45.13"
$ws.Range("E9").Value = "  -0.76%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07463"
$ws.Range("E10").Value = "  -0.45%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.126"
$ws.Range("E11").Value = "  -1.22%  "

$ws.Range("B12").Value = "BinanceUSD"
$ws.Range("C12").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.013"
$ws.Range("E12").Value = "  +1.02%  "

$ws.Range("B13").Value = "Solana"
$ws.Range("C13").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "22.11"
$ws.Range("E13").Value = "  -1.28%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.208"
$ws.Range("E14").Value = "  -0.17%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.342"
$ws.Range("E15").Value = "  +1.35%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.779.83"
$ws.Range("E16").Value = "  -0.67%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001078"
$ws.Range("E17").Value = "  -0.39%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06643"
$ws.Range("E18").Value = "  -0.94%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "81.95"
$ws.Range("E19").Value = "  +0.97%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.010"
$ws.Range("E20").Value = "  +1.08%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.09"
$ws.Range("E21").Value = "  -0.89%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.209"
$ws.Range("E22").Value = "  -2.63%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "27.951.51"
$ws.Range("E23").Value = "  -0.99%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.37"
$ws.Range("E24").Value = "  -4.18%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.421"
$ws.Range("E25").Value = "  +2.01%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "20.38"
$ws.Range("E26").Value = "  -0.05%  "

$ws.Range("B27").Value = "Monero"
$ws.Range("C27").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "151.86"
$ws.Range("E27").Value = "  -1.09%  "

$ws.Range("B28").Value = "LidoDAOToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.347"
$ws.Range("E28").Value = "  -1.63%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.988.60"
$ws.Range("E29").Value = "  -0.39%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.288"
$ws.Range("E30").Value = "  +1.61%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "129.23"
$ws.Range("E31").Value = "  -2.34%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.983"
$ws.Range("E32").Value = "  -2.30%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.773"
$ws.Range("E33").Value = "  -1.63%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.09171"
$ws.Range("E34").Value = "  -2.75%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.2210"
$ws.Range("E35").Value = "  +2.62%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "12.06"
$ws.Range("E36").Value = "  -0.06%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.6617"
$ws.Range("E37").Value = "  -0.09%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.06217"
$ws.Range("E38").Value = "  -0.81%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02296"
$ws.Range("E39").Value = "  -2.89%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.162"
$ws.Range("E40").Value = "  +0.03%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.199"
$ws.Range("E41").Value = "  -1.21%  "

$ws.Range("B42").Value = "WEMIXTOKEN"
$ws.Range("C42").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.437"
$ws.Range("E42").Value = "  -3.03%  "

$ws.Range("B43").Value = "FraxShare"
$ws.Range("C43").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.059"
$ws.Range("E43").Value = "  +0.22%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.009"
$ws.Range("E44").Value = "  +0.97%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "14.00"
$ws.Range("E45").Value = "  +0.32%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.6008"
$ws.Range("E46").Value = "  -1.04%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.811"
$ws.Range("E47").Value = "  -1.48%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "126.14"
$ws.Range("E48").Value = "  -1.73%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.994"
$ws.Range("E49").Value = "  -1.34%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06940"
$ws.Range("E50").Value = "  -2.17%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.124"
$ws.Range("E51").Value = "  -3.48%  "
